$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 200, shifting existing rows 200-286 down to 201-287.
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the new data record.
$ws.Cells.Item(200, 1).Value = 9
$ws.Cells.Item(200, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(200, 3).Value = "Metropolitana"
$ws.Cells.Item(200, 4).Value = 44510
$ws.Cells.Item(200, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(200, 5).Value = 13
$ws.Cells.Item(200, 6).Value = 100112031
$ws.Cells.Item(200, 7).Value = "Poroto verde"
$ws.Cells.Item(200, 8).Value = "Magnum"
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 25
$ws.Cells.Item(200, 11).Value = 33000
$ws.Cells.Item(200, 12).Value = 36000
$ws.Cells.Item(200, 13).Value = 34440
$ws.Cells.Item(200, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(200, 15).Value = "Perú"
$ws.Cells.Item(200, 16).Value = 1378
$ws.Cells.Item(200, 17).Value = 25
$ws.Cells.Item(200, 18).Value = "Hortaliza"
